# The "Mes" column (column C, table "Tabla1") held the month as a plain
# number (1-12). The update replaces those numbers with the Spanish month
# abbreviation used elsewhere in the report ("Ene.", "Feb.", ... "Dic."),
# adding the twelve new labels to the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

# Data rows of Tabla1 run from row 6 to row 85; column C is "Mes".
for ($r = 6; $r -le 85; $r++) {
    $cell = $ws.Range("C$r")
    $monthNum = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNum]
}
